$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 13
$ws.Range("H13").Value = 7466.6665
$ws.Range("I13").Value = 2900
$ws.Range("J13").Value = 9750
$ws.Range("K13").Value = 2900
$ws.Range("L13").Value = 9750
$ws.Range("M13").Value = -2731
$ws.Range("N13").Value = -10088

# Row 31
$ws.Range("H31").Value = 450
$ws.Range("I31").Value = 450
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1350
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1120

# Row 129
$ws.Range("H129").Value = 1120.3889
$ws.Range("I129").Value = 565
$ws.Range("J129").Value = 1231.4667
$ws.Range("K129").Value = 1695
$ws.Range("L129").Value = 3694.4001
$ws.Range("M129").Value = 3305
$ws.Range("N129").Value = -13694.4001

$ws = $wb.Worksheets.Item("ARM")
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# Row 26
$ws.Range("H26").Value = 3236.25
$ws.Range("I26").Value = 1648.3334
$ws.Range("J26").Value = 8000
$ws.Range("K26").Value = 1648.3334
$ws.Range("L26").Value = 8000
$ws.Range("M26").Value = -1318.3334
$ws.Range("N26").Value = -8660

# Row 121
$ws.Range("H121").Value = 37509.168
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 37509.168
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 37509.168
$ws.Range("N121").Value = -41003.168

$ws = $wb.Worksheets.Item("BSM")
# Row 96
$ws.Range("H96").Value = 15500
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 29000
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 29000
$ws.Range("M96").Value = 746
$ws.Range("N96").Value = -34492

# Row 105
$ws.Range("H105").Value = 25436.666
$ws.Range("I105").Value = 36468.332
$ws.Range("J105").Value = 3373.3333
$ws.Range("K105").Value = 36468.332
$ws.Range("L105").Value = 3373.3333
$ws.Range("M105").Value = -34721.332
$ws.Range("N105").Value = -6867.3333

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 43302.5
$ws.Range("I13").Value = 200
$ws.Range("J13").Value = 57670
$ws.Range("K13").Value = 200
$ws.Range("L13").Value = 57670
$ws.Range("M13").Value = -61
$ws.Range("N13").Value = -57948

# Row 22
$ws.Range("H22").Value = 805
$ws.Range("I22").Value = 658.75
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 658.75
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -308.75
$ws.Range("N22").Value = -1700

# Row 25
$ws.Range("H25").Value = 1000000000
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 1000000000
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 1000000000
$ws.Range("N25").Value = -1000000348

# Row 32
$ws.Range("H32").Value = 4433
$ws.Range("I32").Value = 1149.5
$ws.Range("J32").Value = 11000
$ws.Range("K32").Value = 1149.5
$ws.Range("L32").Value = 11000
$ws.Range("M32").Value = -833.5
$ws.Range("N32").Value = -11632

# Row 99
$ws.Range("H99").Value = 15243.637
$ws.Range("I99").Value = 31740
$ws.Range("J99").Value = 1496.6666
$ws.Range("K99").Value = 31740
$ws.Range("L99").Value = 1496.6666
$ws.Range("M99").Value = -30242
$ws.Range("N99").Value = -4492.6666

# Row 122
$ws.Range("H122").Value = 3090614.8
$ws.Range("I122").Value = 13889488
$ws.Range("J122").Value = 5222.2856
$ws.Range("K122").Value = 41668464
$ws.Range("L122").Value = 15666.8568
$ws.Range("M122").Value = -41666014
$ws.Range("N122").Value = -20566.8568

# Row 126
$ws.Range("H126").Value = 15243.637
$ws.Range("I126").Value = 31740
$ws.Range("J126").Value = 1496.6666
$ws.Range("K126").Value = 95220
$ws.Range("L126").Value = 4489.9998
$ws.Range("M126").Value = -92750
$ws.Range("N126").Value = -9429.9998

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 222905.75
$ws.Range("I5").Value = 328.27777
$ws.Range("J5").Value = 371290.75
$ws.Range("K5").Value = 984.83331
$ws.Range("L5").Value = 1113872.25
$ws.Range("M5").Value = -872.83331
$ws.Range("N5").Value = -1114096.25

# Row 31
$ws.Range("H31").Value = 2320
$ws.Range("I31").Value = 500
$ws.Range("J31").Value = 2775
$ws.Range("K31").Value = 1500
$ws.Range("L31").Value = 8325
$ws.Range("M31").Value = -1212
$ws.Range("N31").Value = -8901

# Row 32
$ws.Range("H32").Value = 2166.6667
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2166.6667
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 6500.000100000001
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -7066.000100000001

# Row 135
$ws.Range("H135").Value = 222905.75
$ws.Range("I135").Value = 328.27777
$ws.Range("J135").Value = 371290.75
$ws.Range("K135").Value = 2954.49993
$ws.Range("L135").Value = 3341616.75
$ws.Range("M135").Value = -419.4999299999999
$ws.Range("N135").Value = -3346686.75

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 577055.56
$ws.Range("I12").Value = 577055.56
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 577055.56
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -576915.56
$ws.Range("N12").ClearContents()

# Row 102
$ws.Range("H102").Value = 531341.4
$ws.Range("I102").Value = 1211759.1
$ws.Range("J102").Value = 2127.5557
$ws.Range("K102").Value = 1211759.1
$ws.Range("L102").Value = 2127.5557
$ws.Range("M102").Value = -1210137.1
$ws.Range("N102").Value = -5371.5557

# Row 132
$ws.Range("H132").Value = 3052.8572
$ws.Range("I132").Value = 3547.647
$ws.Range("J132").Value = 2585.5557
$ws.Range("K132").Value = 10642.941
$ws.Range("L132").Value = 7756.6671
$ws.Range("M132").Value = -8112.940999999999
$ws.Range("N132").Value = -12816.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 93
$ws.Range("H93").Value = 681.8182
$ws.Range("I93").Value = 500
$ws.Range("J93").Value = 900
$ws.Range("K93").Value = 500
$ws.Range("L93").Value = 900
$ws.Range("M93").Value = 748

# Row 136
$ws.Range("H136").Value = 4925.015
$ws.Range("I136").Value = 3725.4888
$ws.Range("J136").Value = 7495.4287
$ws.Range("K136").Value = 11176.4664
$ws.Range("L136").Value = 22486.2861
$ws.Range("M136").Value = -8626.466400000001
$ws.Range("N136").Value = -27586.2861

$ws = $wb.Worksheets.Item("WVR")
# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("N31").ClearContents()

# Row 122
$ws.Range("H122").Value = 1284.0769
$ws.Range("I122").Value = 973.75
$ws.Range("J122").Value = 1780.6
$ws.Range("K122").Value = 2921.25
$ws.Range("L122").Value = 5341.799999999999
$ws.Range("M122").Value = -471.25
$ws.Range("N122").Value = -10241.8

# Row 126
$ws.Range("H126").Value = 1055.3
$ws.Range("I126").Value = 758.8333
$ws.Range("J126").Value = 1500
$ws.Range("K126").Value = 2276.4999
$ws.Range("L126").Value = 4500
$ws.Range("M126").Value = 193.5001000000002
$ws.Range("N126").Value = -9440
